# Musa fights Firown (as) Day 98
#
# Row 7 holds references to "blog" card widgets across several cells.
# The rotation advances by one slot:
#   I7: ser 94 -> ser 95
#   E7: ser 95 -> ser 96
#   C7: ser 96 -> ser 98   (new post for "Day 98")
# B7/D7 (other widget types) stay as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 95"
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 96"
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 98"

# Move the active selection to I7 to match the saved view state.
[void]$ws.Range("I7").Select()
